$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Experienced in Python Development." paragraph -> append
# " (Web, Scripting, Automation)" right after "...Development" (before the
# final period), relocate the "_GoBack" bookmark to sit right after that new
# text, and add a single space run right after the bookmark (before the
# period). The pre-existing "Developmen"/"t" run split must be preserved, and
# the newly typed text must land in its own run(s) rather than silently
# merging into its neighbours, so each insertion point is "shielded" with a
# throw-away character-formatting toggle (which forces the engine to keep it
# as a distinct run) before being reverted back to the original formatting
# once the dust has settled (plain property assignment does not trigger a
# run merge, only text insertion does).
# ---------------------------------------------------------------------------

$findRange = $d.Content
$found = $findRange.Find.Execute("Experienced in Python Development.")

$devEnd = $findRange.End - 1          # just after "...Developmen[t]", before "."
$tStart = $devEnd - 1                 # start of the single-letter "t" run

# Shield the trailing "t" run so the upcoming insertion does not fold it back
# into the preceding "Developmen" run.
$tRun = $d.Range($tStart, $devEnd)
$tRun.Font.Bold = 1

$insPoint = $d.Range($devEnd, $devEnd)
$newText = " (Web, Scripting, Automation)"
$insPoint.InsertAfter($newText)
$newTextEnd = $devEnd + $newText.Length

# Split the freshly inserted text back out of the "t" run and restore formatting.
$newTextRange = $d.Range($devEnd, $newTextEnd)
$newTextRange.Font.Bold = 0
$tRunAgain = $d.Range($tStart, $devEnd)
$tRunAgain.Font.Bold = 0

# Move the "_GoBack" bookmark here (it used to sit in the greeting paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($newTextEnd, $newTextEnd))

# Shield the run preceding the insertion point (the text we just added) so the
# new space character becomes its own run instead of merging back into it.
$leftShield = $d.Range($devEnd, $newTextEnd)
$leftShield.Font.Bold = 1

$spaceInsPoint = $d.Range($newTextEnd, $newTextEnd)
$spaceInsPoint.InsertAfter(" ")
$spaceEnd = $newTextEnd + 1

# Differentiate the space from its shielded neighbour using a second,
# independent toggle, then restore both back to plain formatting.
$spaceRange = $d.Range($newTextEnd, $spaceEnd)
$spaceRange.Font.BoldBi = 1

$leftFinal = $d.Range($devEnd, $newTextEnd)
$leftFinal.Font.Bold = 0

$spaceFinal = $d.Range($newTextEnd, $spaceEnd)
$spaceFinal.Font.Bold = 0
$spaceFinal.Font.BoldBi = 0

# ---------------------------------------------------------------------------
# Change 2: remove the stray "s" run (and its now-relocated bookmark) that
# used to sit alone in the paragraph right before "PROFESSIONAL PROFILE",
# leaving that paragraph empty.
# ---------------------------------------------------------------------------

$sFindRange = $d.Content
$sFound = $sFindRange.Find.Execute("s" + [char]13 + "PROFESSIONAL PROFILE")
if (-not $sFound) {
    $sFindRange = $d.Content
    $sFound = $sFindRange.Find.Execute("s")
}
$sStart = $sFindRange.Start
$sRun = $d.Range($sStart, $sStart + 1)
$sRun.Delete()

# ---------------------------------------------------------------------------
# Change 3: merge the "ONSITE - 2" / "0" table-cell runs into one
# "ONSITE - 20" run.
# ---------------------------------------------------------------------------

$dash = [char]8211
$onsiteRange = $d.Content
$onsiteRange.Find.Execute("ONSITE " + $dash + " 20", $true, $false, $false, $false, $false, $true, 1, $false, "ONSITE " + $dash + " 20", 2)
